$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 281, shifting the existing rows 281-294
# down to 283-296 (matches the target diff exactly).
$ws.Rows("281:282").Insert()

# New row 281: Primera, $/caja 12 unidades, updated weekly price data.
$ws.Range("A281").Value = 7
$ws.Range("B281").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C281").Value = "Ñuble"
$ws.Range("D281").Value = 45008
$ws.Range("E281").Value = 16
$ws.Range("F281").Value = "Fruta"
$ws.Range("G281").Value = 100108
$ws.Range("H281").Value = "Tropicales y subtropicales"
$ws.Range("I281").Value = 100108005
$ws.Range("J281").Value = "Piña"
$ws.Range("K281").Value = "Caramelo"
$ws.Range("L281").Value = "Primera"
$ws.Range("M281").Value = 60
$ws.Range("N281").Value = 23000
$ws.Range("O281").Value = 23000
$ws.Range("P281").Value = 23000
$ws.Range("Q281").Value = "$/caja 12 unidades"
$ws.Range("R281").Value = "Ecuador"
$ws.Range("S281").Value = 1917
$ws.Range("T281").Value = 12

# New row 282: Segunda, $/caja 14 unidades, updated weekly price data.
$ws.Range("A282").Value = 7
$ws.Range("B282").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C282").Value = "Ñuble"
$ws.Range("D282").Value = 45008
$ws.Range("E282").Value = 16
$ws.Range("F282").Value = "Fruta"
$ws.Range("G282").Value = 100108
$ws.Range("H282").Value = "Tropicales y subtropicales"
$ws.Range("I282").Value = 100108005
$ws.Range("J282").Value = "Piña"
$ws.Range("K282").Value = "Caramelo"
$ws.Range("L282").Value = "Segunda"
$ws.Range("M282").Value = 80
$ws.Range("N282").Value = 23000
$ws.Range("O282").Value = 23000
$ws.Range("P282").Value = 23000
$ws.Range("Q282").Value = "$/caja 14 unidades"
$ws.Range("R282").Value = "Ecuador"
$ws.Range("S282").Value = 1643
$ws.Range("T282").Value = 14
